$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, M, N, O, P, S

$ws.Range("D2").Value2 = 44209
$ws.Range("M2").Value2 = 100
$ws.Range("N2").Value2 = 10000
$ws.Range("O2").Value2 = 11000
$ws.Range("P2").Value2 = 10500
$ws.Range("S2").Value2 = 750

$ws.Range("D3").Value2 = 44217
$ws.Range("M3").Value2 = 200
$ws.Range("N3").Value2 = 11000
$ws.Range("O3").Value2 = 12000
$ws.Range("P3").Value2 = 11500
$ws.Range("S3").Value2 = 821
